# This script applies a cyclic permutation of data across rows 48-53
# (columns A, B, D, E, F, G, H, P, Q, R) on the "Artfynd" worksheet,
# reproducing the target state described by the diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target values for rows 48..53 (columns A,B,D,E,F,G,H,P,Q,R)
$rows = @{
    48 = @{ A = 111974125; B = 90660; D = "NT"; E = 4362; F = "Blå taggsvamp"; G = "Hydnellum caeruleum"; H = "(Hornem.) P.Karst."; P = "Aloppmoarna i S, Jmt"; Q = 439278.8711310769; R = 6952206.909989387 }
    49 = @{ A = 111974029; B = 88032; D = "VU"; E = 6276; F = "Goliatmusseron"; G = "Tricholoma matsutake"; H = "(S.Ito & S.Imai) Singer"; P = "Aloppmoarna, Jmt"; Q = 439334.7866423383; R = 6952296.802153576 }
    50 = @{ A = 111974124; B = 90666; D = "LC"; E = 4364; F = "Dropptaggsvamp"; G = "Hydnellum ferrugineum"; H = "(Fr.:Fr.) P. Karst."; P = "Aloppmoarna i S, Jmt"; Q = 439276.3867801811; R = 6952196.853249942 }
    51 = @{ A = 111974126; B = 88032; D = "VU"; E = 6276; F = "Goliatmusseron"; G = "Tricholoma matsutake"; H = "(S.Ito & S.Imai) Singer"; P = "Aloppmoarna i S, Jmt"; Q = 439289.9461055733; R = 6952209.002200785 }
    52 = @{ A = 111974134; B = 90658; D = "NT"; E = 4361; F = "Orange taggsvamp"; G = "Hydnellum aurantiacum"; H = "(Batsch:Fr.) P.Karst."; P = "Aloppmoarna i S, Jmt"; Q = 439399.8222122483; R = 6952207.441512506 }
    53 = @{ A = 111974133; B = 90682; D = "NT"; E = 2059; F = "Skrovlig taggsvamp"; G = "Hydnellum scabrosum"; H = "(Fr.) E.Larss., K.H.Larss. & Kõljalg"; P = "Aloppmoarna i S, Jmt"; Q = 439389.9449806474; R = 6952220.480550999 }
}

foreach ($r in $rows.Keys) {
    $data = $rows[$r]
    $ws.Range("A$r").Value = $data.A
    $ws.Range("B$r").Value = $data.B
    $ws.Range("D$r").Value = $data.D
    $ws.Range("E$r").Value = $data.E
    $ws.Range("F$r").Value = $data.F
    $ws.Range("G$r").Value = $data.G
    $ws.Range("H$r").Value = $data.H
    $ws.Range("P$r").Value = $data.P
    $ws.Range("Q$r").Value = $data.Q
    $ws.Range("R$r").Value = $data.R
}
